# part B number 1 complete

$wb = $excel.ActiveWorkbook

# --- Part A: move the stored selection to D5 without disturbing the active tab ---
$wsA = $wb.Worksheets.Item("Part A")
$wsA.Range("D5").Select() | Out-Null

# --- Add the new "Part 2" sheet at the end of the workbook ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Part 2"

# Row 1-2: headers (string insertion order matters for shared-string table order)
$ws3.Range("A1").Value = "Impedance"
$ws3.Range("A2").Value = "Phase a"
$ws3.Range("R2").Value = "Phase c"
$ws3.Range("B1").Value = "Rural Bus 11"
$ws3.Range("J2").Value = "Phase b"

# Row 3: fault-type sub headers
$ws3.Range("B3").Value = "3 Phase"
$ws3.Range("F3").Value = "SLG"
$ws3.Range("J3").Value = "LL"
$ws3.Range("N3").Value = "DLG"
$ws3.Range("R3").Value = "LL"
$ws3.Range("V3").Value = "DLG"

# Row 4: Real / Imag / Mag column headers for each of the six blocks
$ws3.Range("B4").Value = "Real"
$ws3.Range("C4").Value = "Imag"
$ws3.Range("D4").Value = "Mag"

$ws3.Range("F4").Value = "Real"
$ws3.Range("G4").Value = "Imag"
$ws3.Range("H4").Value = "Mag"

$ws3.Range("J4").Value = "Real"
$ws3.Range("K4").Value = "Imag"
$ws3.Range("L4").Value = "Mag"

$ws3.Range("N4").Value = "Real"
$ws3.Range("O4").Value = "Imag"
$ws3.Range("P4").Value = "Mag"

$ws3.Range("R4").Value = "Real"
$ws3.Range("S4").Value = "Imag"
$ws3.Range("T4").Value = "Mag"

$ws3.Range("V4").Value = "Real"
$ws3.Range("W4").Value = "Imag"
$ws3.Range("X4").Value = "Mag"

# Row 5: data values + Mag formulas
$ws3.Range("B5").Value = 0.0183
$ws3.Range("C5").Value = 0.1033
$ws3.Range("D5").Formula = "=SQRT((B5)^2+(C5)^2)"

$ws3.Range("F5").Value = 0.0137
$ws3.Range("G5").Value = 0.1677
$ws3.Range("H5").Formula = "=SQRT((F5)^2+(G5)^2)"

$ws3.Range("J5").Value = -0.7716
$ws3.Range("K5").Value = 0.4124
$ws3.Range("L5").Formula = "=SQRT((J5)^2+(K5)^2)"

$ws3.Range("N5").Value = -0.716
$ws3.Range("O5").Value = 0.3818
$ws3.Range("P5").Formula = "=SQRT((N5)^2+(O5)^2)"

$ws3.Range("R5").Value = -0.7716
$ws3.Range("S5").Value = 0.4124
$ws3.Range("T5").Formula = "=SQRT((R5)^2+(S5)^2)"

$ws3.Range("V5").Value = -0.716
$ws3.Range("W5").Value = 0.3818
$ws3.Range("X5").Formula = "=SQRT((V5)^2+(W5)^2)"

# Row 7-8: zone legend
$ws3.Range("A7").Value = "Impedance"
$ws3.Range("A8").Value = "Zone 1"
$ws3.Range("B8").Value = "Zone 2"
$ws3.Range("C8").Value = "Zone 3"

# Column widths to fit content
$ws3.Columns.Item(1).AutoFit() | Out-Null
$ws3.Columns.Item(2).AutoFit() | Out-Null

# Make the new sheet the active tab/selection, as last step
$ws3.Activate() | Out-Null
$ws3.Range("A9").Select() | Out-Null
